$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2-5 from serial date 45204 to 45207
# (i.e. 2023-10-05 -> 2023-10-08), keeping existing cell formatting/style.
$ws.Range("C2:C5").Value = 45207
